$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SALARY")
$ws.Range("D39").Value = "                              "
